$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Aktiv block ---
$ws.Range("C17").Value = 8
$ws.Range("G17").Value = 4
$ws.Range("C18").Value = "4"
$ws.Range("D19").Value = "6"
$ws.Range("G19").Value = 3
$ws.Range("D20").Value = "37.5% der Karten"
$ws.Range("B21").Value = 5

# --- Karten ohne Aktivitaet ---
$ws.Range("C27").Value = 334
$ws.Range("C28").Value = 21
$ws.Range("B29").Value = "Gruppentreffen 29.06.2020 💩"
$ws.Range("C29").Value = 16
$ws.Range("B30").Value = "Gruppentreffen 06.07.2020"
$ws.Range("C30").Value = 12
$ws.Range("B31").Value = "Gruppentreffen 13.07.2020"
$ws.Range("C31").Value = 1

# --- Gemeinschaftlich: Aktivste Mitglieder ---
$ws.Range("C35").Value = 8
$ws.Range("C36").Value = 7
$ws.Range("C37").Value = 7
$ws.Range("C38").Value = 7
$ws.Range("C39").Value = 7
$ws.Range("F39").Value = 8
$ws.Range("F40").Value = 6
$ws.Range("F41").Value = 6
$ws.Range("F42").Value = 6
$ws.Range("F43").Value = 6

# --- Zuverlaessig: Meiste abgeschlossene Karten ---
$ws.Range("C48").Value = 3
$ws.Range("C49").Value = 3
$ws.Range("C50").Value = 3
$ws.Range("E50").Value = "Maria Lütticke"
$ws.Range("F50").Value = 1
$ws.Range("C51").Value = 3
$ws.Range("E51").Value = "Peter Augustin"
$ws.Range("F51").Value = 1
$ws.Range("C52").Value = 3
$ws.Range("E52").Value = "Noah Brechmann"
$ws.Range("F52").Value = 1
$ws.Range("E53").Value = "Christoph Netsch"
$ws.Range("F53").Value = 1
$ws.Range("E54").Value = "Jacob Escherich"
$ws.Range("F54").Value = 1

# --- Puenktlich ---
$ws.Range("G60").Value = 3

# --- Detailliert ---
$ws.Range("B70").Value = 8
$ws.Range("F70").Value = "Christoph Netsch"
$ws.Range("G70").Value = 2
$ws.Range("B71").Value = 2
$ws.Range("F71").Value = "Eric Pomp"
$ws.Range("G71").Value = 2
$ws.Range("F72").Value = "Jacob Escherich"
$ws.Range("G72").Value = 2
$ws.Range("B73").Value = 1
$ws.Range("F73").Value = "Maria Lütticke"
$ws.Range("G73").Value = 2
$ws.Range("F74").Value = "Marie-Sophie Braun"
$ws.Range("G74").Value = 2
